# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.107.94'
$ws.Range("E2").Value = '  +0.05%  '

$ws.Range("D3").Value = '1.781.60'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.36'
$ws.Range("E5").Value = '  -0.94%  '

$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("E9").Value = '  -1.52%  '

$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("E11").Value = '  +0.81%  '

$ws.Range("D12").Value = '2.038.63'
$ws.Range("E12").Value = '  -0.46%  '

$ws.Range("D14").Value = '1.771.19'
$ws.Range("E14").Value = '  -1.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.623'
$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").Value = '34.089.48'

$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.59'
$ws.Range("E18").Value = '  -0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.85'
$ws.Range("E19").Value = '  +1.07%  '

$ws.Range("E20").Value = '  +1.60%  '

$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.86'
$ws.Range("E22").Value = '  +0.76%  '

$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("E24").Value = '  -0.83%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.92'
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.26'
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("E29").Value = '  +0.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.23'
$ws.Range("E30").Value = '  -1.13%  '

$ws.Range("E31").Value = '  -0.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.71'
$ws.Range("E32").Value = '  +1.38%  '

$ws.Range("E33").Value = '  +2.29%  '

$ws.Range("E34").Value = '  -2.65%  '

$ws.Range("D35").Value = '1.449.95'

$ws.Range("E36").Value = '  +4.89%  '

$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("E39").Value = '  -0.81%  '

$ws.Range("E40").Value = '  +1.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '80.55'
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.72'
$ws.Range("E42").Value = '  +1.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.915'
$ws.Range("E43").Value = '  -0.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.71'
$ws.Range("E44").Value = '  +1.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0520'
$ws.Range("E45").Value = '  +2.41%  '

$ws.Range("E46").Value = '  -0.33%  '

$ws.Range("E47").Value = '  +0.26%  '

$ws.Range("D48").Value = '1.938.22'
$ws.Range("E48").Value = '  -0.58%  '

$ws.Range("E49").Value = '  -5.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.64'
$ws.Range("E50").Value = '  -2.85%  '

$ws.Range("E51").Value = '  +0.20%  '
